$d = $word.ActiveDocument

# Locate the paragraph that ends with the "Nao apenas..." SCRUM sentence, so the
# insertion point is robust even if paragraph indices shift.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*validação das rápidas entregas.*") {
        $target = $para
    }
}
if ($target -eq $null) {
    throw "Could not locate anchor paragraph"
}

$insertPos = $target.Range.End
$rng = $d.Range($insertPos, $insertPos)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>3.</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Organização da equipe </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Por tratar-se de um projeto baseado na metodologia SCRUM, nada mais plausível e eficaz do que seguir o padrão utilizado no mercado. Dentro do time teremos:</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:lastRenderedPageBreak/><w:t>Product</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t>Owner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">– Responsável por administrar o </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Backlog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> do produto, ou seja, irá </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>admistrar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> a demanda do projeto estimando prazos e entregas com todo o time. Tem o papel de alinhar as expectativas do projeto com o solicitante e com o time como um todo. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t>Scrum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> Master</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> – Responsável por retirar os possíveis impedimentos do time de desenvolvimento, e ainda auxiliar na utilização da metodologia e processos do SCRUM.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="pt-BR"/></w:rPr><w:t>Time de Desenvolvimento</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> – Conjunto de profissionais de áreas complementares ao projeto que serão responsáveis pelo desenvolvimento (programação) do projeto em si. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Serão: Programadores (Front e Back </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>end</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>); Cientista de Dados; UX Designer; Analistas para testes;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

$d.Save()
